# Regenerate save_data column G ("K" - strikeouts) values from the source
# box-score data (previously populated from "Strike#", now derived from K).
# Only column G (the "K" column) values change; all other columns are
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 2
    6  = 1
    7  = 1
    8  = 0
    9  = 0
    10 = 2
    11 = 3
    12 = 5
    13 = 3
    14 = 0
    15 = 1
    16 = 2
    17 = 2
    18 = 0
    19 = 2
    20 = 0
    21 = 0
    22 = 2
    23 = 2
    24 = 3
    25 = 0
    26 = 2
    27 = 1
    28 = 0
    29 = 1
    30 = 1
    31 = 4
    32 = 5
    33 = 1
    34 = 1
    36 = 1
    37 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
